$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.306.48'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.868.92'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4702'
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2872'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06584'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08021'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.08'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').Value = '1.871.49'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.117'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6856'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '268.75'
$ws.Range('E16').Value = '  -3.76%  '
$ws.Range('D17').Value = '30.332.78'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.97'
$ws.Range('E18').Value = '  +3.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007659'
$ws.Range('E19').Value = '  +4.97%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '2.117.21'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.281'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.209'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.415'
$ws.Range('E25').Value = '  +2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.59'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.90'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09871'
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.375'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.069'
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04699'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01873'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.43'
$ws.Range('E41').Value = '  -2.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.959'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8414'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4167'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.01'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.212'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.066'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '918.67'
$ws.Range('E49').Value = '  -5.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.53'
$ws.Range('E50').Value = '  +1.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05683'
$ws.Range('E51').Value = '  +0.81%  '
